$wb = $excel.ActiveWorkbook

$wsCreate = $wb.Worksheets.Item("Create Order")
$wsCancel = $wb.Worksheets.Item("Cancel Order")

# --- Cancel Order sheet: fill in the previously-empty "result" column (D) ---
# D2 / D6 / D8 -> "unknown transaction" error payload
# D3           -> successful cancel payload
# D9           -> "can't cancel completed transaction" error payload
#
# Values are entered in the same order the unique strings were first authored
# (so newly-created shared-string entries land at the same indices as the
# target workbook: WALLET/CANCELED payload -> 67, completed-transaction
# error -> 68, unknown-transaction error -> 69).
$wsCancel.Range("D3").Value = '"method":"WALLET","phoneNumber":"081252930398","catalog":{"id":13,"provider":{"id":2,"name":"Telkomsel","image":"https://res.cloudinary.com/alvark/image/upload/v1592209103/danapulsa/Telkomsel_Logo_eviigt_nbbrjv.png"},"value":15000,"price":15000},"status":"CANCELED"'
$wsCancel.Range("D9").Value = '"code":400,"message":"can''t cancel completed transaction"'
$wsCancel.Range("D2").Value = '"code":404,"message":"unknown transaction"'

# D6 previously had the default/general number format (no border-text style); align it
# with the rest of the column (Text format) before writing the value.
$wsCancel.Range("D6").NumberFormat = "@"
$wsCancel.Range("D6").Value = '"code":404,"message":"unknown transaction"'

$wsCancel.Range("D8").Value = '"code":404,"message":"unknown transaction"'

# Widen column D now that it holds much longer JSON payloads.
$wsCancel.Columns.Item(4).ColumnWidth = 241.333333333333

# Move the cell cursor on Cancel Order to D9 (matches where editing finished).
$wsCancel.Range("D9").Select()

# Finish up on the "Create Order" tab/cell, which is the active sheet/selection
# at save time.
$wsCreate.Activate()
$wsCreate.Range("B27").Select()
